# Weekly update: insert a new (most-recent) price observation at row 300
# for "Feria Lagunitas de Puerto Montt - Perejil", pushing all existing
# rows 300..360 down by one (to 301..361) and updating the workbook's
# used-range dimension accordingly (Excel does this automatically when a
# row is inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 300. Excel shifts rows 300-360
# down to 301-361 and copies formatting (incl. the date NumberFormat on
# column D) from the surrounding rows automatically.
$ws.Rows(300).Insert()

# Populate the newly inserted row 300 with the new weekly observation.
# Columns A, B, C, E, F, G, H, I, R are identical across the whole
# market/category block, so copy them straight from the row below
# (the old row 300, now shifted to row 301).
$ws.Cells.Item(300, 1).Value  = $ws.Cells.Item(301, 1).Value2
$ws.Cells.Item(300, 2).Value  = $ws.Cells.Item(301, 2).Value2
$ws.Cells.Item(300, 3).Value  = $ws.Cells.Item(301, 3).Value2
$ws.Cells.Item(300, 4).Value  = 44995
$ws.Cells.Item(300, 5).Value  = $ws.Cells.Item(301, 5).Value2
$ws.Cells.Item(300, 6).Value  = $ws.Cells.Item(301, 6).Value2
$ws.Cells.Item(300, 7).Value  = $ws.Cells.Item(301, 7).Value2
$ws.Cells.Item(300, 8).Value  = $ws.Cells.Item(301, 8).Value2
$ws.Cells.Item(300, 9).Value  = $ws.Cells.Item(301, 9).Value2
$ws.Cells.Item(300, 10).Value = 180
$ws.Cells.Item(300, 11).Value = 7000
$ws.Cells.Item(300, 12).Value = 7000
$ws.Cells.Item(300, 13).Value = 7000
$ws.Cells.Item(300, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(300, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(300, 16).Value = 2333
$ws.Cells.Item(300, 17).Value = 3
$ws.Cells.Item(300, 18).Value = $ws.Cells.Item(301, 18).Value2
